$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

Set-TextValue $ws.Range('D2') '285.55'
Set-TextValue $ws.Range('E2') '1.01%'
Set-TextValue $ws.Range('D3') '29.38'
Set-TextValue $ws.Range('E3') '3.92%'
Set-TextValue $ws.Range('D4') '5.064'
Set-TextValue $ws.Range('E4') '0.88%'
Set-TextValue $ws.Range('D5') '0.06711'
Set-TextValue $ws.Range('E5') '3.22%'
Set-TextValue $ws.Range('D6') '7.315'
Set-TextValue $ws.Range('E6') '1.03%'
Set-TextValue $ws.Range('B7') 'FTXToken'
Set-TextValue $ws.Range('C7') 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D7') '1.381'
Set-TextValue $ws.Range('E7') '-5.55%'
Set-TextValue $ws.Range('B8') 'MXToken'
Set-TextValue $ws.Range('C8') 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D8') '0.9011'
Set-TextValue $ws.Range('E8') '-1.01%'
Set-TextValue $ws.Range('B9') 'WazirX'
Set-TextValue $ws.Range('C9') 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D9') '0.1569'
Set-TextValue $ws.Range('E9') '1.67%'
Set-TextValue $ws.Range('B10') 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws.Range('C10') 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D10') '0.07104'
Set-TextValue $ws.Range('E10') '12.87%'
Set-TextValue $ws.Range('B11') 'MandalaExchangeToken'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D11') '0.07626'
Set-TextValue $ws.Range('E11') '1.33%'
Set-TextValue $ws.Range('B12') 'BitrueCoin'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D12') '0.02919'
Set-TextValue $ws.Range('E12') '4.28%'
Set-TextValue $ws.Range('B13') 'BitMartToken'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D13') '0.08990'
Set-TextValue $ws.Range('E13') '0.29%'
Set-TextValue $ws.Range('B14') 'BitForexToken'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D14') '0.001578'
Set-TextValue $ws.Range('E14') '-0.40%'
Set-TextValue $ws.Range('B15') 'CoinExToken'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D15') '0.04487'
Set-TextValue $ws.Range('E15') '1.42%'
Set-TextValue $ws.Range('B16') 'One'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range('D16') '0.0006477'
Set-TextValue $ws.Range('E16') '1.26%'
Set-TextValue $ws.Range('B17') 'TigerCash'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D17') '0.006219'
Set-TextValue $ws.Range('E17') '0.43%'
Set-TextValue $ws.Range('B18') 'LEO'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D18') '3.449'
Set-TextValue $ws.Range('E18') '0.15%'
Set-TextValue $ws.Range('B19') 'GateToken'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D19') '3.440'
Set-TextValue $ws.Range('E19') '1.82%'
Set-TextValue $ws.Range('D20') '2.228'
Set-TextValue $ws.Range('E20') '-0.53%'
Set-TextValue $ws.Range('D21') '0.3233'
Set-TextValue $ws.Range('E21') '1.13%'
Set-TextValue $ws.Range('D22') '0.1319'
Set-TextValue $ws.Range('E22') '1.06%'
Set-TextValue $ws.Range('D23') '3.935'
Set-TextValue $ws.Range('E23') '-1.22%'
Set-TextValue $ws.Range('D24') '0.1559'
Set-TextValue $ws.Range('E24') '3.29%'
Set-TextValue $ws.Range('D25') '0.001202'
Set-TextValue $ws.Range('E25') '1.67%'
Set-TextValue $ws.Range('D26') '0.004368'
Set-TextValue $ws.Range('E26') '-1.39%'
Set-TextValue $ws.Range('D27') '0.0001170'
Set-TextValue $ws.Range('E27') '-6.15%'
Set-TextValue $ws.Range('D28') '0.0001618'
Set-TextValue $ws.Range('E28') '0.00%'
Set-TextValue $ws.Range('D40') '0.04237'
Set-TextValue $ws.Range('E40') '2.61%'
Set-TextValue $ws.Range('D41') '0.006770'
Set-TextValue $ws.Range('E41') '1.13%'
Set-TextValue $ws.Range('E42') '0.70%'
Set-TextValue $ws.Range('D43') '0.002231'
Set-TextValue $ws.Range('E43') '3.05%'
Set-TextValue $ws.Range('D44') '0.01269'
Set-TextValue $ws.Range('E44') '5.16%'
Set-TextValue $ws.Range('D45') '0.00005756'
Set-TextValue $ws.Range('E45') '4.08%'
Set-TextValue $ws.Range('E46') '0.04%'
Set-TextValue $ws.Range('E47') '15.49%'
